$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 2's formatting down to row 3 (keeps the date style on column A)
$ws.Range("A2:N2").Copy($ws.Range("A3:N3"))

# Fill in the new row's values
$ws.Range("A3").Value = 42606.8828587963
$ws.Range("B3").Value = 26
$ws.Range("C3").Value = 64
$ws.Range("D3").Value = 34
$ws.Range("E3").Value = 50
$ws.Range("F3").Value = 50
$ws.Range("G3").Value = 4663
$ws.Range("H3").Value = 2944
$ws.Range("I3").Value = 164
$ws.Range("J3").Value = 47
$ws.Range("K3").Value = 25
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = "Named"
